$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dish / recipe text fields (shared strings content) ---
# Breakfast hot dish name + recipe number
$ws.Range("C4").Value = "54-19к-2020"
$ws.Range("D4").Value = "Суп молочный с макаронными изделиями"

# Breakfast hot drink recipe number + name
$ws.Range("C5").Value = "54-9гн-2020"
$ws.Range("D5").Value = "Кофейный напиток с молоком"

# Bread row recipe number + name (row 7)
$ws.Range("C7").Value = "54-1з-2020"
$ws.Range("D7").Value = "Сыр твердых сортов"

# Row 8 recipe number + name
$ws.Range("C8").Value = "53-19з2020"
$ws.Range("D8").Value = "Масло сливочное"

# Breakfast 2 / fruit row (row 9) - now filled in
$ws.Range("D9").Value = "Яблоко"

# --- Nutritional numeric values ---
# Row 4
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = 14.36
$ws.Range("G4").Value = 144.7
$ws.Range("H4").Value = 5.26
$ws.Range("I4").Value = 5.52
$ws.Range("J4").Value = 18.4

# Row 5
$ws.Range("F5").Value = 19.6
$ws.Range("G5").Value = 90.8
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 11.1

# Row 6
$ws.Range("G6").Value = 57.9
$ws.Range("H6").Value = 2.3
$ws.Range("I6").Value = 0.3
$ws.Range("J6").Value = 11.5

# Row 7
$ws.Range("E7").Value = 30
$ws.Range("F7").Value = 23.7
$ws.Range("G7").Value = 109.1
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 0

# Row 8
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 8.05
$ws.Range("G8").Value = 71.8
$ws.Range("H8").Value = 0.1
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 17.8

# Row 9 (previously blank E9:J9)
$ws.Range("E9").Value = 150
$ws.Range("F9").Value = 22.5
$ws.Range("G9").Value = 60.6
$ws.Range("H9").Value = 0.6
$ws.Range("I9").Value = 0.6
$ws.Range("J9").Value = 13.5

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(9).RowHeight = 13.8

# --- Selection / active cell ---
$ws.Range("F9").Select()
